# Add new columns I ("I0") and J ("IF") with data for rows 2-31.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - match the formatting already used by the other header cells
# (bold font, thin box border, centered/top aligned).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data values for column I and J, rows 2..31
$iVals = @(8,3,7,9,7,1,7,4,8,6,8,6,7,5,6,5,7,7,4,9,9,6,6,8,7,5,6,8,7,7)
$jVals = @(8,4,8,9,7,1,7,4,8,7,8,7,7,5,6,6,8,8,5,9,9,6,7,9,7,7,7,8,7,7)

for ($r = 0; $r -lt $iVals.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$r]
    $ws.Cells.Item($row, 10).Value = $jVals[$r]
}
